$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDD")

# --- Column I: widened to fit the new review notes ---
$ws.Columns.Item(9).ColumnWidth = 45.3

# --- New note in I16:I18 (merged, centered) - set before I7 so shared-string order matches ---
$ws.Range("I16:I18").HorizontalAlignment = -4108
$ws.Range("I16:I18").VerticalAlignment = -4108
$ws.Range("I16:I18").Merge()
$ws.Range("I16").Value = "refer to rubric! I'm dev + research!"

# --- Row 7 (Constraints) now flagged "Bad"; row 10 (Assumptions) cleared back to normal ---
$ws.Range("D7").Style = "Bad"
$ws.Range("D7").HorizontalAlignment = -4108
$ws.Range("D10").Style = "Normal"
$ws.Range("D10").HorizontalAlignment = -4108

# --- New reviewer note on row 7 ---
$ws.Range("I7").Value = "the real constraint is use of TC & SP score"
$ws.Range("I7").HorizontalAlignment = -4108

# --- Updated word counts / draft completion percentages ---
$ws.Range("E8").Value = 108
$ws.Range("E10").Value = 152
$ws.Range("F10").Value = 0.7
$ws.Range("F11").Value = 0.3
$ws.Range("F12").Value = 0.7
$ws.Range("F14").Value = 0.75

# --- Update the selected cell on the sheet ---
$ws.Range("E21").Select()

$excel.Calculate()
